$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45182 (2023-09-13)
# to 45184 (2023-09-15) for every data row (rows 2 through 452).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 452 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45184
